$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style/formatting used by the
# other header cells (e.g. G1: bold, bordered, centered header style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add corresponding data value in H2 (plain numeric cell, no special style).
$ws.Range("H2").Value = 0
